# TCHD Covid-19 time series update: append row 17 (day 16, date serial 43922)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry down the per-column formatting from row 16 onto the new row 17
# (B = date format, D:F/I:J/K/M = the "bordered/wrapped" numeric format)
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("D16:F16").Copy()
$ws.Range("D17:F17").PasteSpecial(-4122)

$ws.Range("I16:K16").Copy()
$ws.Range("I17:K17").PasteSpecial(-4122)

$ws.Range("M16").Copy()
$ws.Range("M17").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# A17: running day counter, continues the A-column "+1" pattern
$ws.Range("A17").Formula = "=A16+1"

# B17: continuation of the date-increment series
$ws.Range("B17").Formula = "=B16+1"

# C17..M17: new data values for this row
$ws.Range("C17").Value = 242
$ws.Range("D17").Value = 80
$ws.Range("E17").Value = 1158
$ws.Range("F17").Value = 1480
$ws.Range("G17").Value = "#N/A"
$ws.Range("H17").Value = "#N/A"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1
$ws.Range("K17").Formula = "=K16+L17"
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 33

# Restore the author's saved cursor position
$ws.Range("L18").Select() | Out-Null
